# Requirements matrix: add Test Case / Result / Comments / Fixed / Regression Result
# columns (G:K, plus L for an extra comment on row 44), update the "Run command" text
# in C61 (and the consequent shift of the Quit/Exit rows), widen a few columns, resize
# row 61 and update the current selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column C text updates (rows shifted because the old "2 conditions" text was removed) ---
$ws.Range("C61").Value = 'Run command will display a message in 3 conditions.  1. Accepted, including # of transitions to discover.  2. Rejected, including # of transitions to discover. 3. Instantaneous description '
$ws.Range("C62").Value = 'Quit will be invoked by ''q'' or ''Q''.  This will quit any running TM but not exit program.  If not processing string an appropriate message will be displayed.  User is returned to command.'
$ws.Range("C63").Value = 'Exit will be invoked by ''x'' or ''X''. Terminate the application.  If changes were made to input string list output file.  Display message of success or fail to write input string list .str'

# --- New columns G (Test Case #), H (Result), I (Comments), J (Fixed), K (Regression Result), L ---
$ws.Range("G1").Value = 'Test Case'
$ws.Range("H1").Value = 'Result'
$ws.Range("I1").Value = 'Comments'
$ws.Range("J1").Value = 'Fixed'
$ws.Range("K1").Value = 'Regression Result'
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("G18").Value = 1
$ws.Range("G19").Value = 1
$ws.Range("G20").Value = 1
$ws.Range("G21").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 'pass'
$ws.Range("I24").Value = 'Error message displayed'
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 'pass'
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 'fail'
$ws.Range("I26").Value = 'empty string was represented as \'
$ws.Range("J26").Value = 'x'
$ws.Range("K26").Value = 'pass'
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 'fail'
$ws.Range("I27").Value = 'empty string was ignored.'
$ws.Range("J27").Value = 'x'
$ws.Range("K27").Value = 'pass'
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 'pass'
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 'pass'
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 'fail'
$ws.Range("J30").Value = 'x'
$ws.Range("K30").Value = 'pass'
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 'fail'
$ws.Range("J31").Value = 'x'
$ws.Range("K31").Value = 'pass'
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 'pass'
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 'fail'
$ws.Range("I33").Value = 'File could not be opened/saved.'
$ws.Range("J33").Value = 'x'
$ws.Range("K33").Value = 'pass'
$ws.Range("G34").Value = 'x'
$ws.Range("G35").Value = 'x'
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 'pass'
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 'pass'
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 'pass'
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 'pass'
$ws.Range("G40").Value = 6
$ws.Range("H40").Value = 'PASS'
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 'pass'
$ws.Range("G42").Value = 8
$ws.Range("H42").Value = 'pass'
$ws.Range("G43").Value = 9
$ws.Range("H43").Value = 'pass'
$ws.Range("G44").Value = 10
$ws.Range("H44").Value = 'fail'
$ws.Range("I44").Value = 'help is displayed for main menu, but not individual commands.'
$ws.Range("J44").Value = 'x'
$ws.Range("K44").Value = 'pass'
$ws.Range("L44").Value = 'Help added to List function and removed from being displayed at main menu.'
$ws.Range("G45").Value = 11
$ws.Range("H45").Value = 'pass'
$ws.Range("G46").Value = 12
$ws.Range("H46").Value = 'pass'
$ws.Range("G47").Value = 13
$ws.Range("H47").Value = 'pass'
$ws.Range("G48").Value = 14
$ws.Range("H48").Value = 'pass'
$ws.Range("G49").Value = 15
$ws.Range("H49").Value = 'pass'
$ws.Range("G50").Value = 16
$ws.Range("H50").Value = 'pass'
$ws.Range("G51").Value = 17
$ws.Range("H51").Value = 'pass'
$ws.Range("G52").Value = 18
$ws.Range("H52").Value = 'pass'
$ws.Range("G53").Value = 19
$ws.Range("H53").Value = 'pass'
$ws.Range("G54").Value = 20
$ws.Range("H54").Value = 'pass'
$ws.Range("G55").Value = 21
$ws.Range("H55").Value = 'pass'
$ws.Range("G56").Value = 22
$ws.Range("H56").Value = 'pass'
$ws.Range("G57").Value = 23
$ws.Range("H57").Value = 'pass'
$ws.Range("G58").Value = 24
$ws.Range("H58").Value = 'pass'
$ws.Range("G59").Value = 25
$ws.Range("H59").Value = 'pass'
$ws.Range("G60").Value = 26
$ws.Range("H60").Value = 'pass'
$ws.Range("G61").Value = 27
$ws.Range("H61").Value = 'pass'
$ws.Range("G62").Value = 28
$ws.Range("H62").Value = 'pass'
$ws.Range("G63").Value = 29
$ws.Range("H63").Value = 'fail'
$ws.Range("I63").Value = 'exitting is hard'
$ws.Range("J63").Value = 'x'
$ws.Range("K63").Value = 'pass'

# --- Formatting: bold header cells, wrap the one G cell that mirrors the bestFit column C style ---
$ws.Range("G1","H1","I1","J1","K1").Font.Bold = $true
$ws.Range("G49").WrapText = $true

# --- Column widths for the new columns ---
$ws.Columns.Item(7).ColumnWidth = 18.9453125
$ws.Columns.Item(8).ColumnWidth = 16.015625
$ws.Columns.Item(9).ColumnWidth = 12.6953125

# --- Row height + view/selection state ---
$ws.Rows.Item(61).RowHeight = 75
$ws.Range("I63").Select()
